# MVPlayerOTY.xlsx update:
#  - Header: A1 "Year" -> "Season" (B1 stays "Player")
#  - Column A: "Year" (e.g. 2017) -> "Season" range string (e.g. 2017-18)
#  - Column B: "Player Name" -> "Player Name\playerid"
#  - Column O ("Award"/"MVP") removed entirely
#  - Selection moved to O2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row Season string + Player\playerid string (team/stat columns C:N are unchanged).
$rows = @(
    @{R=2; S="2017-18"; P="James Harden\hardeja01"},
    @{R=3; S="2016-17"; P="Russell Westbrook\westbru01"},
    @{R=4; S="2015-16"; P="Stephen Curry\curryst01"},
    @{R=5; S="2014-15"; P="Stephen Curry\curryst01"},
    @{R=6; S="2013-14"; P="Kevin Durant\duranke01"},
    @{R=7; S="2012-13"; P="LeBron James\jamesle01"},
    @{R=8; S="2011-12"; P="LeBron James\jamesle01"},
    @{R=9; S="2010-11"; P="Derrick Rose\rosede01"},
    @{R=10; S="2009-10"; P="LeBron James\jamesle01"},
    @{R=11; S="2008-09"; P="LeBron James\jamesle01"},
    @{R=12; S="2007-08"; P="Kobe Bryant\bryanko01"},
    @{R=13; S="2006-07"; P="Dirk Nowitzki\nowitdi01"},
    @{R=14; S="2005-06"; P="Steve Nash\nashst01"},
    @{R=15; S="2004-05"; P="Steve Nash\nashst01"},
    @{R=16; S="2003-04"; P="Kevin Garnett\garneke01"},
    @{R=17; S="2002-03"; P="Tim Duncan\duncati01"},
    @{R=18; S="2001-02"; P="Tim Duncan\duncati01"},
    @{R=19; S="2000-01"; P="Allen Iverson\iversal01"},
    @{R=20; S="1999-00"; P="Shaquille O'Neal\onealsh01"},
    @{R=21; S="1998-99"; P="Karl Malone\malonka01"},
    @{R=22; S="1997-98"; P="Michael Jordan\jordami01"},
    @{R=23; S="1996-97"; P="Karl Malone\malonka01"},
    @{R=24; S="1995-96"; P="Michael Jordan\jordami01"},
    @{R=25; S="1994-95"; P="David Robinson\robinda01"},
    @{R=26; S="1993-94"; P="Hakeem Olajuwon\olajuha01"},
    @{R=27; S="1992-93"; P="Charles Barkley\barklch01"},
    @{R=28; S="1991-92"; P="Michael Jordan\jordami01"},
    @{R=29; S="1990-91"; P="Michael Jordan\jordami01"},
    @{R=30; S="1989-90"; P="Magic Johnson\johnsma02"},
    @{R=31; S="1988-89"; P="Magic Johnson\johnsma02"},
    @{R=32; S="1987-88"; P="Michael Jordan\jordami01"},
    @{R=33; S="1986-87"; P="Magic Johnson\johnsma02"},
    @{R=34; S="1985-86"; P="Larry Bird\birdla01"},
    @{R=35; S="1984-85"; P="Larry Bird\birdla01"},
    @{R=36; S="1983-84"; P="Larry Bird\birdla01"},
    @{R=37; S="1982-83"; P="Moses Malone\malonmo01"},
    @{R=38; S="1981-82"; P="Moses Malone\malonmo01"},
    @{R=39; S="1980-81"; P="Julius Erving\ervinju01"},
    @{R=40; S="1979-80"; P="Kareem Abdul-Jabbar\abdulka01"}
)

# Header row.
$ws.Range("A1").Value = "Season"
$ws.Range("B1").Value = "Player"

foreach ($row in $rows) {
    $ws.Cells.Item($row.R, 1).Value = $row.S
    $ws.Cells.Item($row.R, 2).Value = $row.P
}

# Drop the Award column (O) entirely.
$ws.Columns("O:O").Delete()

# Match the saved selection from the authored workbook.
$ws.Range("O2").Select()
